$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 19.6
$ws.Range("B3").Value = 4.7
$ws.Range("C3").Value = 37.8
$ws.Range("C4").Value = 40.5
$ws.Range("C5").Value = 34.4
$ws.Range("C6").Value = 17.6
$ws.Range("C7").Value = 45.9
$ws.Range("C9").Value = 24.7
$ws.Range("C11").Value = 13.9
$ws.Range("C13").Value = 15.4
$ws.Range("C14").Value = 108.6
$ws.Range("C16").Value = 115.1
$ws.Range("C20").Value = 9.800000000000001
$ws.Range("C21").Value = 2.7
